$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-08-13 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-08-14 Thursday", 2)

# Update the division problems in the table. Cells are addressed by
# (row, column) rather than by text search, since several of the new
# values collide with other original/new values elsewhere in the table
# (e.g. "56÷5=" is both a pre-existing cell value and the replacement
# for another cell), which would make a simple global Find/Replace
# ambiguous or order-dependent.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "55÷9="
$t.Cell(1, 2).Range.Text = "27÷2="
$t.Cell(1, 3).Range.Text = "62÷3="
$t.Cell(1, 4).Range.Text = "16÷5="
$t.Cell(1, 5).Range.Text = "65÷2="

$t.Cell(5, 1).Range.Text = "63÷5="
$t.Cell(5, 2).Range.Text = "56÷5="
$t.Cell(5, 3).Range.Text = "27÷2="
$t.Cell(5, 4).Range.Text = "29÷3="
$t.Cell(5, 5).Range.Text = "79÷3="

$t.Cell(9, 1).Range.Text = "53÷3="
$t.Cell(9, 2).Range.Text = "84÷7="
$t.Cell(9, 3).Range.Text = "51÷4="
$t.Cell(9, 4).Range.Text = "10÷9="
$t.Cell(9, 5).Range.Text = "42÷9="

$t.Cell(13, 1).Range.Text = "15÷4="
$t.Cell(13, 2).Range.Text = "73÷8="
$t.Cell(13, 3).Range.Text = "23÷9="
$t.Cell(13, 4).Range.Text = "18÷5="
$t.Cell(13, 5).Range.Text = "68÷7="

$t.Cell(17, 1).Range.Text = "35÷9="
$t.Cell(17, 2).Range.Text = "59÷8="
$t.Cell(17, 3).Range.Text = "85÷7="
$t.Cell(17, 4).Range.Text = "47÷2="
$t.Cell(17, 5).Range.Text = "24÷2="
